# Automatische test-sync: 2025-08-01 23:41:50
# Append the new "Testmail #3" log row to the "Logs" sheet, grow the
# conditional-formatting ranges to cover it, and bump the "Overig" tally
# on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 8

$ws.Cells.Item($row, 1).Value = "Kun jij dit afhandelen?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #3: Kun jij dit afhandelen?"
$ws.Cells.Item($row, 4).Value = "Overig"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-01 23:40:52"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# Grow each conditional-formatting block so it covers the new row too.
$ws.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D8"))
$ws.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G8"))
$ws.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H8"))
$ws.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I8"))
$ws.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J8"))

# Dashboard: bump the "Overig" count from 5 to 6.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
